$d = $word.ActiveDocument
$vt = [char]11

# ---------------------------------------------------------------------------
# 1. "<  com.philips...CircleIndicator" paragraph: merge the " " run and the
#    "com.philips..." run into a single run (no visible text change).
#    A temporary zero-length bookmark anchors the run boundary right before
#    the replacement so the preceding "<" run is left untouched.
# ---------------------------------------------------------------------------
$pCircle = $d.Paragraphs.Item(8)
$rCircle = $pCircle.Range
$circleText = $rCircle.Text
$spaceIdx = $circleText.IndexOf(" com.philips")
if ($spaceIdx -ge 0) {
    $spanStart = $rCircle.Start + $spaceIdx
    $spanEnd = $rCircle.End

    $rAnchor = $d.Range($spanStart, $spanStart)
    $d.Bookmarks.Add("ztempAnchor1", $rAnchor) | Out-Null

    $rSpan = $d.Range($spanStart, $spanEnd)
    $rSpan.Text = " com.philips.cdp.uikit.dotnavigation.CircleIndicator"

    if ($d.Bookmarks.Exists("ztempAnchor1")) {
        $d.Bookmarks.Item("ztempAnchor1").Delete()
    }
}

# ---------------------------------------------------------------------------
# 2. "android:layout_height / android:layout_width" paragraph:
#    - add 3 leading spaces before "android:layout_height..."
#    - change the 4 leading spaces before "android:layout_width..." to 3
#      and drop the trailing ' />'
#    - a new paragraph is then added right after it for
#      'android:padding="10dp" />'
# ---------------------------------------------------------------------------
$pLayout = $d.Paragraphs.Item(9)
$rLayout = $pLayout.Range
$layoutStart = $rLayout.Start
$layoutEnd = $rLayout.End
$layoutText = $rLayout.Text
$breakIdx = $layoutText.IndexOf($vt)

# Replace the part after the manual line break first (so offsets for the
# first part remain valid).
$rAfterBreak = $d.Range($layoutStart + $breakIdx + 1, $layoutEnd - 1)
$rAfterBreak.Text = '   android:layout_width="wrap_content"'

$rBeforeBreak = $d.Range($layoutStart, $layoutStart + $breakIdx)
$rBeforeBreak.Text = '   android:layout_height="wrap_content"'

# Insert the new paragraph right after this one.
$pLayout = $d.Paragraphs.Item(9)
$pLayout.Range.InsertParagraphAfter() | Out-Null
$pPadding = $d.Paragraphs.Item(10)
$pPadding.Alignment = 0
$pPadding.Range.Text = '                             android:padding="10dp" />'

# ---------------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark from the "To Set UnSelected..." paragraph
#    to right before 'android:layout_width="wrap_content"' in the paragraph
#    just edited above.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks.Item("_GoBack")
    $bm.Delete()
}

$pLayout = $d.Paragraphs.Item(9)
$rLayout = $pLayout.Range
$layoutText2 = $rLayout.Text
$breakIdx2 = $layoutText2.IndexOf($vt)
$bookmarkPos = $rLayout.Start + $breakIdx2 + 1 + 3
$rPoint = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $rPoint) | Out-Null

# ---------------------------------------------------------------------------
# 4. "To Set UnSelected Circle width and height ..." paragraph: merge the
#    "To Set " / "Un" / "Selected Circle ... using " runs into a single run,
#    and merge the " / set" / "Un" / "SelectedCircleHeight" runs into a
#    single run (no visible text change - this is where the "_GoBack"
#    bookmark used to live before step 3 moved it away).
# ---------------------------------------------------------------------------
$pUnSel = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("To Set UnSelected Circle")) {
        $pUnSel = $cand
        break
    }
}

if ($pUnSel -ne $null) {
    $rUnSel = $pUnSel.Range
    $unSelText = $rUnSel.Text
    $unSelStart = $rUnSel.Start

    # -- merge "To Set " + "Un" + "Selected Circle ... using " --
    $part1Marker = "using "
    $part1End = $unSelText.IndexOf($part1Marker) + $part1Marker.Length
    $spanStart1 = $unSelStart
    $spanEnd1 = $unSelStart + $part1End

    $rAnchor1 = $d.Range($spanEnd1, $spanEnd1)
    $d.Bookmarks.Add("ztempAnchor2", $rAnchor1) | Out-Null

    $rSpan1 = $d.Range($spanStart1, $spanEnd1)
    $rSpan1.Text = "To Set UnSelected Circle width and height either set attributes or programmatically using "

    if ($d.Bookmarks.Exists("ztempAnchor2")) {
        $d.Bookmarks.Item("ztempAnchor2").Delete()
    }

    # -- merge " / set" + "Un" + "SelectedCircleHeight" --
    $rUnSel = $pUnSel.Range
    $unSelText = $rUnSel.Text
    $unSelStart = $rUnSel.Start
    $part2Marker = " / set"
    $spanStart2 = $unSelStart + $unSelText.IndexOf($part2Marker)
    $spanEnd2 = $rUnSel.End

    $rAnchor2 = $d.Range($spanStart2, $spanStart2)
    $d.Bookmarks.Add("ztempAnchor3", $rAnchor2) | Out-Null

    $rSpan2 = $d.Range($spanStart2, $spanEnd2)
    $rSpan2.Text = " / setUnSelectedCircleHeight"

    if ($d.Bookmarks.Exists("ztempAnchor3")) {
        $d.Bookmarks.Item("ztempAnchor3").Delete()
    }
}

Write-Output "Done"
